$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Simple Level-of-Effort value edits ---
$ws.Range("F11").Value = 20
$ws.Range("F16").Value = 40

# --- Insert a new blank separator row above row 28 ---
# (This shifts the existing rows 28-37 down to 29-38, inheriting the
#  formatting of the row above, which already matches the blank
#  separator row's style.)
$ws.Rows(28).Insert()

# After the insert, the three rows that used to be 28/29/30 (Release = 2)
# are now rows 29/30/31; mark their Release column with "?" instead of 2.
$ws.Range("A29").Value = "?"
$ws.Range("A30").Value = "?"
$ws.Range("A31").Value = "?"

# Restore the active cell selection as recorded in the saved file.
$ws.Range("H18").Select()
